$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.112.98"
$ws.Range("E2").Value = "  -1.64%  "

$ws.Range("D3").Value = "2.174.70"
$ws.Range("E3").Value = "  -1.93%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "'250.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.18%  "

$ws.Range("D6").Value = "'0.610"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.15%  "

$ws.Range("D7").Value = "'66.14"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -7.33%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").Value = "'0.587"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.22%  "

$ws.Range("D10").Value = "'58.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.77%  "

$ws.Range("D11").Value = "'36.26"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -10.53%  "

$ws.Range("D12").Value = "'0.0933"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.04%  "

$ws.Range("D13").Value = "'0.103"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.06%  "

$ws.Range("D14").Value = "'6.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.85%  "

$ws.Range("D15").Value = "2.503.37"
$ws.Range("E15").Value = "  -1.94%  "

$ws.Range("D16").Value = "'14.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.68%  "

$ws.Range("D17").Value = "'0.842"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.19%  "

$ws.Range("D18").Value = "2.178.25"
$ws.Range("E18").Value = "  -2.34%  "

$ws.Range("D19").Value = "41.061.58"
$ws.Range("E19").Value = "  -1.55%  "

$ws.Range("D20").Value = "0.0₃0942"
$ws.Range("E20").Value = "  -2.06%  "

$ws.Range("D21").Value = "'71.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.79%  "

$ws.Range("D22").Value = "'6.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.80%  "

$ws.Range("D23").Value = "'229.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.12%  "

$ws.Range("D24").Value = "'2.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.58%  "

$ws.Range("D25").Value = "'3.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.98%  "

$ws.Range("E26").Value = "  +0.19%  "

$ws.Range("D27").Value = "'11.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.14%  "

$ws.Range("D28").Value = "'2.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.14%  "

$ws.Range("D29").Value = "'167.89"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.58%  "

$ws.Range("D30").Value = "'2.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.55%  "

$ws.Range("D31").Value = "'20.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.74%  "

$ws.Range("E32").Value = "  -1.63%  "

$ws.Range("D33").Value = "'5.63"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.33%  "

$ws.Range("D34").Value = "'0.0744"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.98%  "

$ws.Range("E35").Value = "  -2.16%  "

$ws.Range("D36").Value = "'4.52"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.21%  "

$ws.Range("D37").Value = "'3.96"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.29%  "

$ws.Range("D38").Value = "'24.53"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.98%  "

$ws.Range("D39").Value = "'0.0303"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.21%  "

$ws.Range("D40").Value = "'5.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +12.97%  "

$ws.Range("E41").Value = "  -3.81%  "

$ws.Range("D42").Value = "'5.51"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.17%  "

$ws.Range("D43").Value = "'60.55"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.94%  "

$ws.Range("D44").Value = "'11.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.03%  "

$ws.Range("D45").Value = "'8.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.50%  "

$ws.Range("B46").Value = "BinanceUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D46").Value = "'1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.15%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.0990"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.75%  "

$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "'0.188"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.35%  "

$ws.Range("D49").Value = "'1.13"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.83%  "

$ws.Range("D50").Value = "'4.27"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -9.94%  "

$ws.Range("E51").Value = "  -3.94%  "

